# Update column F (dSF) values for a set of rows, per repull/push of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 1
    4  = 1
    6  = -4
    8  = -8
    9  = -1
    10 = 5
    13 = -4
    14 = 3
    15 = -1
    16 = -1
    18 = 1
    19 = 4
    20 = 1
    21 = -1
    25 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
